$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.726.79'
$ws.Range("E2").Value = '  +2.20%  '
$ws.Range("D3").Value = '1.874.32'
$ws.Range("E3").Value = '  +2.12%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = '326.23'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("D7").Value = '0.4666'
$ws.Range("E7").Value = '  +0.46%  '
$ws.Range("D8").Value = '0.3911'
$ws.Range("E8").Value = '  +0.94%  '
$ws.Range("D9").Value = '0.07892'
$ws.Range("E9").Value = '  -0.09%  '
$ws.Range("D10").Value = '0.9755'
$ws.Range("E10").Value = '  +1.58%  '
$ws.Range("D11").Value = '22.18'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("D12").Value = '1.858.73'
$ws.Range("E12").Value = '  -6.57%  '
$ws.Range("D13").Value = '5.713'
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '6.971'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '0.06955'
$ws.Range("E15").Value = '  +1.99%  '
$ws.Range("D16").Value = '88.35'
$ws.Range("E16").Value = '  +1.68%  '
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '0.00001006'
$ws.Range("E18").Value = '  +0.73%  '
$ws.Range("D19").Value = '16.87'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").Value = '28.773.26'
$ws.Range("E21").Value = '  +2.18%  '
$ws.Range("D22").Value = '5.322'
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = '11.04'
$ws.Range("E23").Value = '  +0.35%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.286.48'
$ws.Range("E24").Value = '  +5.91%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '2.128'
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("D26").Value = '153.11'
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").Value = '19.27'
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").Value = '5.713'
$ws.Range("E28").Value = '  -0.36%  '
$ws.Range("D29").Value = '1.988'
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("D30").Value = '119.54'
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("D31").Value = '0.09364'
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("D32").Value = '0.9294'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("D33").Value = '5.291'
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").Value = '''1.350'
$ws.Range("E34").Value = '  +2.05%  '
$ws.Range("D35").Value = '''3.350'
$ws.Range("E35").Value = '  +0.14%  '
$ws.Range("D36").Value = '0.05856'
$ws.Range("E36").Value = '  -0.84%  '
$ws.Range("D37").Value = '0.02119'
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("D38").Value = '1.146'
$ws.Range("E38").Value = '  +0.26%  '
$ws.Range("D39").Value = '7.832'
$ws.Range("E39").Value = '  +1.16%  '
$ws.Range("D40").Value = '0.5692'
$ws.Range("E40").Value = '  +1.62%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '9.961'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.1786'
$ws.Range("E42").Value = '  +0.79%  '
$ws.Range("D43").Value = '0.07274'
$ws.Range("E43").Value = '  +3.77%  '
$ws.Range("D44").Value = '11.79'
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("D45").Value = '1.177'
$ws.Range("E45").Value = '  -4.52%  '
$ws.Range("D46").Value = '0.5333'
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("D47").Value = '1.836'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").Value = '113.35'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").Value = '2.053'
$ws.Range("E49").Value = '  -5.55%  '
$ws.Range("D50").Value = '''2.370'
$ws.Range("E50").Value = '  +2.13%  '
$ws.Range("D51").Value = '1.005'
$ws.Range("E51").Value = '  +0.59%  '
